# The paragraph before the trailing "_GoBack" bookmark contained a stray
# placeholder run ("asas") left over from drafting. Remove it now that the
# real creative-solution text (brainstorming / "Lluvia de ideas") lives
# elsewhere in the document, leaving the bookmark in an otherwise empty
# paragraph, exactly as in the target revision.

$d = $word.ActiveDocument

$range = $d.Content
$range.Find.Execute("asas", $true, $true, $false, $false, $false, $true, 1, $false, "", 2)
